$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each (cell, new text value) pair below reproduces the refreshed crypto-ticker
# snapshot (Price / Volume(1h) columns) captured by the GitHub Actions job.
# NumberFormat is forced to Text ("@") before the write so Excel stores the
# value as a literal string (matching the original inline-string cells) rather
# than re-interpreting number-looking text like "310.78" or "1.17%" as a number.
$updates = @{
    'D2' = '310.78'
    'E2' = '1.17%'
    'D3' = '37.59'
    'E3' = '0.13%'
    'D4' = '5.116'
    'E4' = '0.23%'
    'D5' = '0.07855'
    'E5' = '-0.27%'
    'D6' = '4.395'
    'E6' = '1.24%'
    'D7' = '1.899'
    'E7' = '-3.79%'
    'D8' = '8.223'
    'E8' = '-0.06%'
    'D9' = '2.863'
    'E9' = '-8.58%'
    'E11' = '-6.68%'
    'D12' = '0.1907'
    'E12' = '0.25%'
    'D13' = '0.09399'
    'E13' = '5.85%'
    'D14' = '0.03417'
    'E14' = '-0.27%'
    'D15' = '0.09613'
    'E15' = '-1.56%'
    'D16' = '0.001361'
    'E16' = '-2.04%'
    'D17' = '0.005829'
    'E17' = '-1.15%'
    'D18' = '3.544'
    'E18' = '-1.20%'
    'D19' = '0.3432'
    'E19' = '-0.06%'
    'D20' = '5.256'
    'E20' = '4.80%'
    'E21' = '-0.73%'
    'D22' = '0.2584'
    'E22' = '3.47%'
    'D23' = '0.02100'
    'E23' = '179.72%'
    'D24' = '0.04346'
    'E24' = '0.56%'
    'D25' = '0.001199'
    'E25' = '-1.75%'
    'D26' = '0.004273'
    'E26' = '-7.05%'
    'D27' = '0.0001298'
    'E27' = '-63.87%'
    'E39' = '-8.72%'
    'D40' = '0.05065'
    'E40' = '1.04%'
    'D41' = '0.007626'
    'E41' = '1.70%'
    'E42' = '-8.17%'
    'D43' = '0.1353'
    'E43' = '-0.09%'
    'D44' = '0.002067'
    'E44' = '2.27%'
    'D45' = '0.008592'
    'E45' = '7.07%'
    'D46' = '0.00006689'
    'E46' = '2.28%'
    'D47' = '0.00000000749'
    'E47' = '-0.39%'
    'E48' = '-0.49%'
    'D49' = '0.002930'
    'E49' = '-2.40%'
    'D50' = '0.00002097'
    'E50' = '-0.39%'
    'D51' = '0.0001997'
    'E51' = '-0.39%'
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cell]
}

